$wb = $excel.ActiveWorkbook

# Rename the sheets: Sheet1 -> Logic, Sheet2 -> Codes
$wsLogic = $wb.Worksheets.Item("Sheet1")
$wsLogic.Name = "Logic"

$wsCodes = $wb.Worksheets.Item("Sheet2")
$wsCodes.Name = "Codes"

# Activate the "Codes" sheet so it becomes the selected/active tab
$wsCodes.Activate()
